$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.036.19"
$ws.Range("E2").Value = "  -2.79%  "
$ws.Range("D3").Value = "2.648.39"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.97"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.32"
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.02"
$ws.Range("E9").Value = "  +8.81%  "
$ws.Range("E10").Value = "  -3.60%  "
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "3.116.46"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").Value = "59.091.53"
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.03"
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").Value = "2.664.44"
$ws.Range("E16").Value = "  -3.54%  "
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "339.98"
$ws.Range("E18").Value = "  -3.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.37"
$ws.Range("E19").Value = "  -4.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.35"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.36"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.09"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.418"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").Value = "0.0₃0801"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.09"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.67"
$ws.Range("E29").Value = "  -3.45%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.84"
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.30"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.15"
$ws.Range("E34").Value = "  -4.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.19"
$ws.Range("E35").Value = "  -3.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.892"
$ws.Range("E36").Value = "  -6.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.873"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.72"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.47"
$ws.Range("E39").Value = "  -6.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.59"
$ws.Range("E40").Value = "  -3.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.616"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.95"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "274.76"
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0969"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.65"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0533"
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.031.88"
$ws.Range("E48").Value = "  -5.16%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.78"
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0228"
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.82"
$ws.Range("E51").Value = "  -2.01%  "
